# Opp to Eng validation and expense req changes
#
# - AddOpportunity: Retainer/MonthlyFee/ContingentFee/Fee values were stored
#   as "10.0" -- normalize to "10". Est. Transaction Size/Market Cap (SICCode
#   column) was "9999.0" -- normalize to "9999".
# - AddContact: the FVA "required fields" validation message (val3/column K)
#   no longer calls out the Valuation Date field, so drop that trailing
#   clause from the message.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AddOpportunity")
$ws1.Range("O2").Value = "10"
$ws1.Range("P2").Value = "10"
$ws1.Range("Q2").Value = "10"
$ws1.Range("T2").Value = "9999"
$ws1.Range("AB2").Value = "10"

$ws3 = $wb.Worksheets.Item("AddContact")
$ws3.Range("K2").Value = "Opportunity Detail - Client: Street Address., Opportunity Detail - Client: City Address., Opportunity Detail - Client: Postal Code., Opportunity Detail - Subject: Street Address., Opportunity Detail - Subject: City Address., Opportunity Detail - Subject: Postal Code Address."
# Shortening the message reflows the wrapped row height (225pt -> 216pt).
$ws3.Rows.Item(2).RowHeight = 216

# Restore the post-edit selection state seen in the commit: AddOpportunity
# cursor left on AB2 (the last cell touched), AddContact re-activated with
# the cursor on J2, and AddContact remains the active tab.
$ws1.Range("AB2").Select()

$ws3.Activate()
$ws3.Range("J2").Select()
